$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: title paragraph ("To-do task 6") -
#   - remove the centered alignment (w:jc)
#   - replace the two original runs with five runs:
#       43 spaces (bold, 18pt)
#       "To-do task 6" (bold, 18pt)
#       22 spaces (bold, 18pt)
#       "Lec: " (bold, 14pt)
#       "ArrayList 3" (italic, underlined, 14pt)
# -----------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Alignment = 0

$spaces43 = "".PadRight(43)
$spaces22 = "".PadRight(22)

$frag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$frag = $frag + '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">' + $spaces43 + '</w:t></w:r>'
$frag = $frag + '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>To-do task 6</w:t></w:r>'
$frag = $frag + '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t xml:space="preserve">' + $spaces22 + '</w:t></w:r>'
$frag = $frag + '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Lec: </w:t></w:r>'
$frag = $frag + '<w:r><w:rPr><w:i/><w:iCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>ArrayList 3</w:t></w:r>'
$frag = $frag + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$titleRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$titleRange.InsertXML($frag)

# -----------------------------------------------------------------------
# Change 2: "Seekbar increasing ... listview." paragraph -
#   merge the three runs (split apart by spell-check proofErr markers
#   around "listview") back into a single plain run.
# -----------------------------------------------------------------------
$target = "Seekbar increasing and making table below in listview."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $pr = $para.Range
    if ($pr.Text -like "*$target*") {
        $full = $d.Range($pr.Start, $pr.End - 1)
        $full.Delete()
        $ins = $d.Range($pr.Start, $pr.Start)
        $ins.InsertAfter($target)
        break
    }
}
